# German translations for HIVE TEAMS.docx
# Applies the set of text replacements described by the commit's diff.
# wdFindContinue = 1, wdReplaceAll = 2 (used throughout)

$d = $word.ActiveDocument
$nbsp = [char]160

function Replace-InParagraph {
    param(
        [int]$ParaIndex,
        [string]$FindText,
        [string]$ReplaceText
    )
    $p = $d.Paragraphs.Item($ParaIndex)
    $rng = $p.Range
    $ok = $rng.Find.Execute($FindText, $false, $false, $false, $false, $false, $true, 1, $false, $ReplaceText, 2)
    if (-not $ok) {
        Write-Output "WARNING: replacement not found in paragraph $ParaIndex for '$FindText'"
    }
}

# --- HIVE TEAM: QUALITÄTSSICHERUNG section ---
Replace-InParagraph 45 "security" "Sicherheit"
Replace-InParagraph 47 "QA Tester" "QS-Prüfer"

# --- HIVE TEAM: DEVELOPMENT section ---
$find50 = $nbsp + "HIVE TEAM: DEVELOPMENT"
Replace-InParagraph 50 $find50 "HIVE TEAM: ENTWICKLUNG"
Replace-InParagraph 51 "Responsible for building SmartCash and supporting applications." "Verantwortlich für den Aufbau von SmartCash und unterstützenden Anwendungen."
Replace-InParagraph 53 "Hive Coordinator" "Hive Koordinator"
Replace-InParagraph 54 "Creator of the Dash N Drink Soda Machine & SmartCash POS." "Urheber von Dash N Drink Soda Machine & SmartCash POS."
Replace-InParagraph 56 "Developer" "Entwickler"
Replace-InParagraph 58 "C++ Software Engineer" "C++ Softwareingenieur"
Replace-InParagraph 60 "Developer" "Entwickler"
Replace-InParagraph 62 "Developer" "Entwickler"

# --- HIVE TEAM: OUTREACH 2 section ---
$find63 = $nbsp + "HIVE TEAM: OUTREACH 2"
Replace-InParagraph 63 $find63 "  HIVE TEAM: ÖFFENTLICHKEITSARBEIT 2"
Replace-InParagraph 64 "This team focuses on community building, growth, general user acquisition in South America" "Dieses Team konzentriert sich auf Gemeinschaftsbildung, Wachstum und allgemeine Nutzerakquise in Südamerika"
Replace-InParagraph 66 "Hive Coordinator" "Hive Koordinator"
Replace-InParagraph 68 "Outreach Support" "Öffentlichkeitsarbeit Support"
Replace-InParagraph 70 "Outreach Support" "Öffentlichkeitsarbeit Support"
Replace-InParagraph 72 "Outreach Support" "Öffentlichkeitsarbeit Support"

# --- HIVE TEAM: SUPPORT / WEB heading paragraph (multiple runs) ---
$find75a = $nbsp + "HIVE TEAM: SUPPORT" + $nbsp
Replace-InParagraph 75 $find75a "HIVE TEAM: SUPPORT "
$find75b = $nbsp + "WEB"
Replace-InParagraph 75 $find75b "WEB"

Write-Output "Done."
